$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New header cells
$ws.Cells.Item(1, 10).Value = "ps"
$ws.Cells.Item(1, 11).Value = "e_ps"
$ws.Cells.Item(1, 12).Value = "Amp"

# New data values for rows 11-19 (columns J, K, L)
$data = @(
    @(11, 0.57, 0.41, 417),
    @(12, 0.34, 0.12, 118),
    @(13, 0.17, 0.27, 269),
    @(14, 0.22, 0.15, 250),
    @(15, 0.33, 0.17, 241),
    @(16, 0.55, 0.2, 231),
    @(17, 0.27, 0.28, 513),
    @(18, 0.27, 0.19, 544),
    @(19, -0.28, 0.22, 281)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 10).Value = $row[1]
    $ws.Cells.Item($r, 11).Value = $row[2]
    $ws.Cells.Item($r, 12).Value = $row[3]
}

# Move the active selection to reflect where the analyst was last working
$ws.Range("K10").Select()

